# Generate Report for Handback
# Adds a new handback-status row (file a2c1d7f3-05cc-4888-95b8-cccdb897c336.md)
# to the Overview sheet and to the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$fileId   = "a2c1d7f3-05cc-4888-95b8-cccdb897c336.md"
$ext      = ".md"
$status   = "Handed back: in sync with en-US"

# -----------------------------------------------------------------
# Sheet "Overview" - adds row 4
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()
$r = $rowOverview.Range.Row

$wsOverview.Cells.Item($r, 1).Value = $fileId
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($r, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5305a1d361bf3adfacbec1b1a3dc5c6e12960dda/e2e/$fileId",
    "",
    "",
    "e2e\$fileId"
)
$wsOverview.Cells.Item($r, 3).Value = $ext
$wsOverview.Cells.Item($r, 5).Value = $status
$wsOverview.Cells.Item($r, 6).Value = $status
$wsOverview.Cells.Item($r, 7).Value = "2016-08-28 02:44:08"
$wsOverview.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -----------------------------------------------------------------
# Sheet "zh-cn" - adds row 4
# -----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $loZhCn.ListRows.Add()
$r = $rowZhCn.Range.Row

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($r, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5305a1d361bf3adfacbec1b1a3dc5c6e12960dda/e2e/$fileId",
    "",
    "",
    $fileId
)
$wsZhCn.Cells.Item($r, 2).Value = $ext
$wsZhCn.Cells.Item($r, 3).Value = $status
$wsZhCn.Cells.Item($r, 4).Value = "e2e"
$wsZhCn.Cells.Item($r, 5).Value = "ht"
$wsZhCn.Cells.Item($r, 6).Value = "'True"
$wsZhCn.Cells.Item($r, 7).Value = "a2c1d7f3-05cc-4888-95b8-cccdb897c336.5305a1d361bf3adfacbec1b1a3dc5c6e12960dda.zh-cn.xlf"
$wsZhCn.Cells.Item($r, 8).Value = "2016-08-28 02:44:00"
$wsZhCn.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($r, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5305a1d361bf3adfacbec1b1a3dc5c6e12960dda/e2e/$fileId",
    "",
    "",
    $fileId
)
$wsZhCn.Cells.Item($r, 10).Value = "a2c1d7f3-05cc-4888-95b8-cccdb897c336.5305a1d361bf3adfacbec1b1a3dc5c6e12960dda.zh-cn.xlf"
$wsZhCn.Cells.Item($r, 11).Value = "2016-08-28 02:44:28"
$wsZhCn.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($r, 12).Value = "'"
$wsZhCn.Cells.Item($r, 13).Value = "'True"
$wsZhCn.Cells.Item($r, 14).Value = "'"
$wsZhCn.Cells.Item($r, 15).Value = "'False"
$wsZhCn.Cells.Item($r, 16).Value = "'"

# -----------------------------------------------------------------
# Sheet "de-de" - adds row 4
# -----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $loDeDe.ListRows.Add()
$r = $rowDeDe.Range.Row

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($r, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5305a1d361bf3adfacbec1b1a3dc5c6e12960dda/e2e/$fileId",
    "",
    "",
    $fileId
)
$wsDeDe.Cells.Item($r, 2).Value = $ext
$wsDeDe.Cells.Item($r, 3).Value = $status
$wsDeDe.Cells.Item($r, 4).Value = "e2e"
$wsDeDe.Cells.Item($r, 5).Value = "ht"
$wsDeDe.Cells.Item($r, 6).Value = "'True"
$wsDeDe.Cells.Item($r, 7).Value = "a2c1d7f3-05cc-4888-95b8-cccdb897c336.5305a1d361bf3adfacbec1b1a3dc5c6e12960dda.de-de.xlf"
$wsDeDe.Cells.Item($r, 8).Value = "2016-08-28 02:44:08"
$wsDeDe.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($r, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5305a1d361bf3adfacbec1b1a3dc5c6e12960dda/e2e/$fileId",
    "",
    "",
    $fileId
)
$wsDeDe.Cells.Item($r, 10).Value = "a2c1d7f3-05cc-4888-95b8-cccdb897c336.5305a1d361bf3adfacbec1b1a3dc5c6e12960dda.de-de.xlf"
$wsDeDe.Cells.Item($r, 11).Value = "2016-08-28 02:44:35"
$wsDeDe.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($r, 12).Value = "'"
$wsDeDe.Cells.Item($r, 13).Value = "'True"
$wsDeDe.Cells.Item($r, 14).Value = "'"
$wsDeDe.Cells.Item($r, 15).Value = "'False"
$wsDeDe.Cells.Item($r, 16).Value = "'"
